# Kapitel 8 eliminiert und alles dahinter nach vorne gezogen
#
# Insert a "Summe" column (B+C per row) right after the existing
# Theorie/Praxis columns, and fold the old stand-alone "Ziel / h" row into
# the total row right above it (which is renamed "Ziel" and now sums from
# the new Summe column instead of re-adding Theorie+Praxis).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before D; old D/E shift right to E/F.
$ws.Columns("D").Insert()
$ws.Columns("D").ColumnWidth = 6

# 2) New "Summe" header and per-row totals (Theorie + Praxis).
$ws.Range("D1").Value = "Summe"
$ws.Range("D2").Formula = "=+B2+C2"
$ws.Range("D3:D12").Formula = "=+B3+C3"
# Row 11 is (and stays) a blank spacer row - drop the fill-down leftover.
$ws.Range("D11").ClearContents()

# Row 8 is part of a custom-formatted block (s=2); the freshly inserted
# cell inherited that row style, but it should stay unstyled like the
# other Summe cells, so re-pull plain formatting from a neighbour.
$ws.Range("D9").Copy()
$ws.Range("D8").PasteSpecial(-4122)

# 3) Row 14 ("Gesamtsumme") absorbs old row 15 ("Ziel"/"h"): rename it,
#    point the grand total at the new Summe column, keep the hours
#    conversion and unit label, then remove the now-empty row 15.
$ws.Range("A14").Value = "Ziel"
$ws.Range("B14").Formula = "=+D12"
$ws.Range("C14").Formula = "=B14/60"
$ws.Range("E14").Value = "h"
$ws.Rows("15").Delete()

# 4) Match the editor's final selection.
$ws.Range("D11").Select()
